$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.81788957118988
$ws.Range("B1").Value = 4.539745330810547
$ws.Range("C1").Value = 4.017266273498535
$ws.Range("D1").Value = 0.9045870900154114
$ws.Range("E1").Value = 0.4760756194591522
